# "Generate Report for Handback"
#
# Updates the handback-status report for the row corresponding to
# e2e\8220be52-0502-477a-8ec6-8b132a9b5049.md after a fresh handback
# cycle: the Overview sheet's "Latest HO Xliff Generate Date", and both
# locale sheets' "Correspond Handoff Datetime" / "Correspond Handback
# DateTime" columns for that file move forward to the new timestamps
# recorded by the handback run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 4 is the 8220be52-...-md entry -----------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2017-02-21 05:15:58"

# --- zh-cn sheet: row 4 is the same file's zh-cn handoff/handback -------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2017-02-21 05:15:41"
$wsZhCn.Range("L4").Value = "2017-02-21 05:17:58"

# --- de-de sheet: row 4 is the same file's de-de handoff/handback -------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2017-02-21 05:15:58"
$wsDeDe.Range("L4").Value = "2017-02-21 05:18:21"
